$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from the row above (D3) into D4 first so the new date
# cell picks up the existing date number-format style instead of Excel
# minting a brand new numFmt entry.
$ws.Range("D3").Copy($ws.Range("D4"))

$ws.Range("A4").Value = "Lammii  Diroo"
$ws.Range("B4").Value = 922956646
$ws.Range("C4").Value = 10000
$ws.Range("D4").Value = 46027

$ws.Range("D6").Select()
